$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose entire data (columns A:AY) need to be swapped.
$pairs = @(
    @(7, 8),
    @(15, 16),
    @(27, 29),
    @(28, 30)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("A" + $r1 + ":AY" + $r1)
    $range2 = $ws.Range("A" + $r2 + ":AY" + $r2)

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
